$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 591.6667
$ws.Range("J17").Value = 591.6667
$ws.Range("L17").Value = 1775.0001
$ws.Range("N17").Value = -2111.0001
$ws.Range("H38").Value = 714.2174
$ws.Range("I38").Value = 134.5
$ws.Range("J38").Value = 1023.4
$ws.Range("K38").Value = 403.5
$ws.Range("L38").Value = 3070.2
$ws.Range("M38").Value = -31.5
$ws.Range("N38").Value = -3814.2
$ws.Range("H116").Value = 1852.6923
$ws.Range("I116").Value = 1609.4445
$ws.Range("J116").Value = 2400
$ws.Range("K116").Value = 1609.4445
$ws.Range("L116").Value = 2400
$ws.Range("M116").Value = 1832.5555
$ws.Range("N116").Value = -9284
$ws.Range("H137").Value = 2274224
$ws.Range("I137").Value = 3572659.8
$ws.Range("J137").Value = 1961.875
$ws.Range("K137").Value = 10717979.4
$ws.Range("L137").Value = 5885.625
$ws.Range("M137").Value = -10715429.4
$ws.Range("N137").Value = -10985.625
$ws.Range("H141").Value = 1325.9778
$ws.Range("I141").Value = 1348
$ws.Range("J141").Value = 852.5
$ws.Range("K141").Value = 4044
$ws.Range("L141").Value = 2557.5
$ws.Range("M141").Value = 1136
$ws.Range("N141").Value = -12917.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14502.202
$ws.Range("I32").Value = 17273.473
$ws.Range("J32").Value = 7112.148
$ws.Range("K32").Value = 17273.473
$ws.Range("L32").Value = 7112.148
$ws.Range("M32").Value = -16986.473
$ws.Range("N32").Value = -7686.148
$ws.Range("H76").Value = 40000
$ws.Range("J76").Value = 40000
$ws.Range("L76").Value = 40000
$ws.Range("N76").Value = -40676
$ws.Range("H79").Value = 40000
$ws.Range("J79").Value = 40000
$ws.Range("L79").Value = 40000
$ws.Range("N79").Value = -42340
$ws.Range("H122").Value = 4832494
$ws.Range("I122").Value = 1647.8182
$ws.Range("J122").Value = 111111110
$ws.Range("K122").Value = 4943.4546
$ws.Range("L122").Value = 333333330
$ws.Range("M122").Value = -2493.4546
$ws.Range("N122").Value = -333338230
$ws.Range("H127").Value = 29931.25
$ws.Range("J127").Value = 29931.25
$ws.Range("L127").Value = 29931.25
$ws.Range("N127").Value = -39851.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1020.58826
$ws.Range("I20").Value = 1012.3077
$ws.Range("K20").Value = 1012.3077
$ws.Range("M20").Value = -765.3077
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H94").Value = 1183.3334
$ws.Range("I94").Value = 550
$ws.Range("J94").Value = 1500
$ws.Range("K94").Value = 550
$ws.Range("L94").Value = 1500
$ws.Range("M94").Value = -99
$ws.Range("N94").Value = -2402
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H128").Value = 1310
$ws.Range("I128").Value = 1310
$ws.Range("K128").Value = 3930
$ws.Range("M128").Value = -1440
$ws.Range("H134").Value = 2152.9688
$ws.Range("I134").Value = 1228.5
$ws.Range("J134").Value = 4926.375
$ws.Range("K134").Value = 3685.5
$ws.Range("L134").Value = 14779.125
$ws.Range("M134").Value = -1150.5
$ws.Range("N134").Value = -19849.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2812.413
$ws.Range("I31").Value = 1056.697
$ws.Range("J31").Value = 7269.231
$ws.Range("K31").Value = 1056.697
$ws.Range("L31").Value = 7269.231
$ws.Range("M31").Value = -761.6969999999999
$ws.Range("N31").Value = -7859.231
$ws.Range("H34").Value = 2812.413
$ws.Range("I34").Value = 1056.697
$ws.Range("J34").Value = 7269.231
$ws.Range("K34").Value = 1056.697
$ws.Range("L34").Value = 7269.231
$ws.Range("M34").Value = -854.6969999999999
$ws.Range("N34").Value = -7673.231
$ws.Range("H58").Value = 16950536
$ws.Range("I58").Value = 19609010
$ws.Range("J58").Value = 2762.625
$ws.Range("K58").Value = 19609010
$ws.Range("L58").Value = 2762.625
$ws.Range("M58").Value = -19608807
$ws.Range("N58").Value = -3168.625
$ws.Range("H132").Value = 85411.25
$ws.Range("I132").Value = 54876.95
$ws.Range("J132").Value = 201441.6
$ws.Range("K132").Value = 164630.85
$ws.Range("L132").Value = 604324.8
$ws.Range("M132").Value = -162100.85
$ws.Range("N132").Value = -609384.8
$ws.Range("H136").Value = 16950536
$ws.Range("I136").Value = 19609010
$ws.Range("J136").Value = 2762.625
$ws.Range("K136").Value = 58827030
$ws.Range("L136").Value = 8287.875
$ws.Range("M136").Value = -58824480
$ws.Range("N136").Value = -13387.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1048.4
$ws.Range("I92").Value = 1083
$ws.Range("J92").Value = 996.5
$ws.Range("K92").Value = 3249
$ws.Range("L92").Value = 2989.5
$ws.Range("M92").Value = -2001
$ws.Range("N92").Value = -5485.5
$ws.Range("H114").Value = 25641738
$ws.Range("I114").Value = 664.1111
$ws.Range("J114").Value = 83334160
$ws.Range("K114").Value = 1992.3333
$ws.Range("L114").Value = 250002480
$ws.Range("M114").Value = 1261.6667
$ws.Range("N114").Value = -250008988
$ws.Range("H117").Value = 4167751
$ws.Range("I117").Value = 519.75
$ws.Range("J117").Value = 8334982.5
$ws.Range("K117").Value = 1559.25
$ws.Range("L117").Value = 25004947.5
$ws.Range("M117").Value = 1882.75
$ws.Range("N117").Value = -25011831.5
$ws.Range("H120").Value = 565225.4399999999
$ws.Range("I120").Value = 1253007.5
$ws.Range("J120").Value = 14999.8
$ws.Range("K120").Value = 3759022.5
$ws.Range("L120").Value = 44999.39999999999
$ws.Range("M120").Value = -3754184.5
$ws.Range("N120").Value = -54675.39999999999
$ws.Range("H121").Value = 42955564
$ws.Range("I121").Value = 1082.375
$ws.Range("J121").Value = 53062500
$ws.Range("K121").Value = 3247.125
$ws.Range("L121").Value = 159187500
$ws.Range("M121").Value = -1937.125
$ws.Range("N121").Value = -159190120
$ws.Range("H133").Value = 5718
$ws.Range("I133").Value = 2789.8572
$ws.Range("K133").Value = 8369.571599999999
$ws.Range("M133").Value = -3309.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 88716.664
$ws.Range("I70").Value = 204420
$ws.Range("J70").Value = 6071.4287
$ws.Range("K70").Value = 204420
$ws.Range("L70").Value = 6071.4287
$ws.Range("M70").Value = -204150
$ws.Range("N70").Value = -6611.4287
$ws.Range("H73").Value = 88716.664
$ws.Range("I73").Value = 204420
$ws.Range("J73").Value = 6071.4287
$ws.Range("K73").Value = 204420
$ws.Range("L73").Value = 6071.4287
$ws.Range("M73").Value = -203484
$ws.Range("N73").Value = -7943.4287
$ws.Range("H107").Value = 288.57895
$ws.Range("I107").Value = 98
$ws.Range("J107").Value = 500.33334
$ws.Range("K107").Value = 98
$ws.Range("L107").Value = 500.33334
$ws.Range("M107").Value = 1822
$ws.Range("N107").Value = -4340.33334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3423.3076
$ws.Range("I40").Value = 3286.1428
$ws.Range("K40").Value = 3286.1428
$ws.Range("M40").Value = -3150.1428
$ws.Range("H100").Value = 1295.2
$ws.Range("I100").Value = 1119
$ws.Range("K100").Value = 1119
$ws.Range("M100").Value = -578
$ws.Range("H115").Value = 60302
$ws.Range("J115").Value = 60302
$ws.Range("L115").Value = 60302
$ws.Range("N115").Value = -62652
$ws.Range("H136").Value = 66085
$ws.Range("I136").Value = 43192.5
$ws.Range("J136").Value = 151114.28
$ws.Range("K136").Value = 129577.5
$ws.Range("L136").Value = 453342.84
$ws.Range("M136").Value = -127027.5
$ws.Range("N136").Value = -458442.84

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 85328.71000000001
$ws.Range("I136").Value = 77791.08
$ws.Range("J136").Value = 94236.82000000001
$ws.Range("K136").Value = 233373.24
$ws.Range("L136").Value = 282710.46
$ws.Range("M136").Value = -287810.46
